$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.172.40"
$ws.Range("E2").Value = "  -0.76%  "

$ws.Range("D3").Value = "'1.630.00"
$ws.Range("E3").Value = "  -1.26%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "'215.92"
$ws.Range("E5").Value = "  -0.82%  "

$ws.Range("D6").Value = "'0.518"
$ws.Range("E6").Value = "  +1.14%  "

$ws.Range("E7").Value = "  +0.40%  "

$ws.Range("E8").Value = "  -0.68%  "

$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").Value = "'20.26"
$ws.Range("E10").Value = "  +0.66%  "

$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "'1.629.00"
$ws.Range("E12").Value = "  -1.17%  "

$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").Value = "'0.543"
$ws.Range("E14").Value = "  +0.12%  "

$ws.Range("D15").Value = "'27.163.99"
$ws.Range("E15").Value = "  -0.56%  "

$ws.Range("D16").Value = "'64.60"
$ws.Range("E16").Value = "  -4.80%  "

$ws.Range("D17").Value = "'0.0₃0733"
$ws.Range("E17").Value = "  -1.07%  "

$ws.Range("D18").Value = "'215.82"
$ws.Range("E18").Value = "  -2.03%  "

$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").Value = "'6.90"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("D22").Value = "'2.48"
$ws.Range("E22").Value = "  -1.14%  "

$ws.Range("D23").Value = "'9.10"
$ws.Range("E23").Value = "  -1.26%  "

$ws.Range("D24").Value = "'147.89"
$ws.Range("E24").Value = "  +0.20%  "

$ws.Range("E25").Value = "  +0.29%  "

$ws.Range("D26").Value = "'7.27"
$ws.Range("E26").Value = "  -3.98%  "

$ws.Range("E27").Value = "  -0.47%  "

$ws.Range("E28").Value = "  -1.64%  "

$ws.Range("E29").Value = "  -0.62%  "

$ws.Range("E30").Value = "  -0.93%  "

$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("E32").Value = "  -1.32%  "

$ws.Range("D33").Value = "'1.312.07"
$ws.Range("E33").Value = "  +3.59%  "

$ws.Range("E34").Value = "  -2.27%  "

$ws.Range("D35").Value = "'2.46"
$ws.Range("E35").Value = "  +0.24%  "

$ws.Range("E36").Value = "  -1.88%  "

$ws.Range("D37").Value = "'0.850"
$ws.Range("E37").Value = "  +0.85%  "

$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("E39").Value = "  +0.26%  "

$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("D41").Value = "'0.801"
$ws.Range("E41").Value = "  -1.29%  "

$ws.Range("D42").Value = "'63.73"
$ws.Range("E42").Value = "  +1.86%  "

$ws.Range("D43").Value = "'1.767.31"
$ws.Range("E43").Value = "  -1.19%  "

$ws.Range("D44").Value = "'5.19"
$ws.Range("E44").Value = "  -4.37%  "

$ws.Range("D45").Value = "'90.66"
$ws.Range("E45").Value = "  -1.34%  "

$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("E47").Value = "  +6.49%  "

$ws.Range("E48").Value = "  +20.56%  "

$ws.Range("D49").Value = "'0.0515"
$ws.Range("E49").Value = "  +0.52%  "

$ws.Range("D50").Value = "'7.55"
$ws.Range("E50").Value = "  -1.85%  "

$ws.Range("D51").Value = "'0.0956"
$ws.Range("E51").Value = "  -2.29%  "
